$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for rows 2-10, columns G,H,I,J,M,N,O,P,Q,R,S,T
$data = @{
    2 = @{ G=26.144619; H=78.433857; I=0.5211737020083955; J=0.5211737020083955;
           M=135.0916853333333; N=405.2750559999999; O=0.7123704212620513; P=0.7123704212620514;
           Q=3531.920643107888; R=31787.28578797099; S=0.3712687296504235; T=0.3712687296504236 }
    3 = @{ G=26.144619; H=78.433857; I=0.5211737020083955; J=0.5211737020083955;
           O=0.2125756143240238; P=0.2125756143240238;
           Q=1053.946343142976; R=9485.517088286784; S=0.1107888198739604; T=0.1107888198739604 }
    4 = @{ G=26.144619; H=78.433857; I=0.5211737020083955; J=0.5211737020083955;
           M=14.23299766666667; N=42.698993; O=0.07505396441392481; P=0.07505396441392483;
           Q=372.116301222889; R=3349.046711006001; S=0.03911615248401158; T=0.03911615248401158 }
    5 = @{ I=0.3571392594830743; J=0.3571392594830742;
           M=135.0916853333333; N=405.2750559999999; O=0.7123704212620513; P=0.7123704212620514;
           Q=2420.2823706792; R=21782.5413361128; S=0.2544154447271747; T=0.2544154447271747 }
    6 = @{ I=0.3571392594830743; J=0.3571392594830742;
           O=0.2125756143240238; P=0.2125756143240238;
           Q=722.2268028384001; R=6500.0412255456; S=0.07591909748384146; T=0.07591909748384144 }
    7 = @{ I=0.3571392594830743; J=0.3571392594830742;
           M=14.23299766666667; N=42.698993; O=0.07505396441392481; P=0.07505396441392483;
           Q=254.9962512463501; R=2294.96626121715; S=0.02680471727205812; T=0.02680471727205812 }
    8 = @{ G=6.104416333333333; H=18.313249; I=0.1216870385085301; J=0.1216870385085301;
           M=135.0916853333333; N=405.2750559999999; O=0.7123704212620513; P=0.7123704212620514;
           Q=824.6558904463269; R=7421.903014016943; S=0.08668624688445307; T=0.08668624688445309 }
    9 = @{ G=6.104416333333333; H=18.313249; I=0.1216870385085301; J=0.1216870385085301;
           O=0.2125756143240238; P=0.2125756143240238;
           Q=246.0822730497209; R=2214.740457447488; S=0.02586769696622193; T=0.02586769696622194 }
    10 = @{ G=6.104416333333333; H=18.313249; I=0.1216870385085301; J=0.1216870385085301;
            M=14.23299766666667; N=42.698993; O=0.07505396441392481; P=0.07505396441392483;
            Q=86.88414342869522; R=781.957290858257; S=0.009133094657855119; T=0.009133094657855121 }
}

foreach ($rowNum in $data.Keys) {
    $cols = $data[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $ws.Range("$colLetter$rowNum").Value = $cols[$colLetter]
    }
}
